# "updated main GSC export data"
#
# The "Chart" sheet (Worksheets(1)) holds a rolling window of one row per
# day: column A is the date, column B is always 0 (Invalid), column C is
# the cumulative "Valid" count. This export rolled the window forward by
# two days: the oldest day (2025-10-16) drops off the top, and two new
# days (2026-01-12, 2026-01-13) are appended at the bottom with the latest
# counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the oldest dated row (row 2 = 2025-10-16, Valid count 50). Deleting
# the row shifts every later row up by one, which is exactly the rest of
# the diff (each remaining row now shows the next day's figures).
$ws.Rows("2:2").Delete()

# Append the two new trailing days. Row 89 used to be the last row
# (2026-01-11); the new rows are 89 and 90.
$ws.Range("A89").Value = "'2026-01-12"
$ws.Range("A89").ClearFormats()
$ws.Range("B89").Value = 0
$ws.Range("C89").Value = 26

$ws.Range("A90").Value = "'2026-01-13"
$ws.Range("A90").ClearFormats()
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 26
